$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 6, shifting existing rows 6-13 down to 7-14.
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 (same pattern as the other "Primera" rows,
# with a new date and volume).
$ws.Range("A6").Value = 7
$ws.Range("B6").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C6").Value = "Ñuble"
$ws.Range("D6").Value = 44664
$ws.Range("D6").Style = $ws.Range("D7").Style
$ws.Range("D6").NumberFormat = $ws.Range("D7").NumberFormat
$ws.Range("E6").Value = 16
$ws.Range("F6").Value = 100112043
$ws.Range("G6").Value = "Pepino dulce"
$ws.Range("H6").Value = "Cultivar IV Región"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 160
$ws.Range("K6").Value = 15000
$ws.Range("L6").Value = 16000
$ws.Range("M6").Value = 15500
$ws.Range("N6").Value = "$/bandeja 18 kilos"
$ws.Range("O6").Value = "Provincia de Limarí"
$ws.Range("P6").Value = 861
$ws.Range("Q6").Value = 18
$ws.Range("R6").Value = "Hortaliza"
